$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (field1 / field2); everything below shifts up one
# row — this also turns the old B1 string "field2" into the numeric 1 that
# used to be in B2, etc., matching the target data layout exactly.
$ws.Rows.Item(1).Delete()

# Narrow column A slightly (was 11.52 "chars", now ~11.34).
$ws.Columns.Item(1).ColumnWidth = 10.42

# Selection moves back to A1 (was B1 before the edit).
$ws.Range("A1").Select() | Out-Null
